# Added results for fine grained classification Config 6 (Label Powerset,
# Thunderbird sheet): F1/Precision/Recall/Accuracy columns (C:F) for rows
# 45-49, 51-55, 57-61 (ngram(3), ngram(2), ngram(1) + POS configurations).
$wb = $excel.ActiveWorkbook
$thunderbird = $wb.Worksheets.Item("Thunderbird")

$thunderbird.Range("C45").Value = "0.602 0.469 0.302 0.310 0.446"
$thunderbird.Range("D45").Value = "0.523 0.508 0.360 0.356 0.535"
$thunderbird.Range("E45").Value = "0.462 0.312 0.179 0.184 0.288"
$thunderbird.Range("F45").Value = "0.768 0.845 0.917 0.915 0.975"

$thunderbird.Range("C46").Value = "0.685 0.590 0.256 0.412 0.460"
$thunderbird.Range("D46").Value = "0.440 0.458 0.264 0.300 0.369"
$thunderbird.Range("E46").Value = "0.636 0.438 0.147 0.263 0.300 "
$thunderbird.Range("F46").Value = "0.716 0.831 0.908 0.900 0.969"

$thunderbird.Range("C47").Value = "0.568 0.643 0.390 0.481 0.489"
$thunderbird.Range("D47").Value = "0.542 0.503 0.259 0.417 0.473"
$thunderbird.Range("E47").Value = "0.418 0.498 0.246 0.320 0.325"
$thunderbird.Range("F47").Value = "0.773 0.844 0.894 0.917 0.973"

$thunderbird.Range("C48").Value = "0.389 0.258 0.044 0.060 0.399"
$thunderbird.Range("D48").Value = "0.389 0.453 0.106 0.226 0.513 "
$thunderbird.Range("E48").Value = "0.250 0.149 0.022 0.031 0.250"
$thunderbird.Range("F48").Value = "0.723 0.838 0.915 0.921 0.974"

$thunderbird.Range("C49").Value = "0.477 0.259 0.179 0.305 0.416"
$thunderbird.Range("D49").Value = "0.588 0.900 0.846 0.820 0.955 "
$thunderbird.Range("E49").Value = "0.321 0.149 0.098 0.180 0.263"
$thunderbird.Range("F49").Value = "0.781 0.864 0.933 0.936 0.981"

$thunderbird.Range("C51").Value = "0.593 0.478 0.269 0.376 0.504"
$thunderbird.Range("D51").Value = "0.527 0.500 0.304 0.449 0.614"
$thunderbird.Range("E51").Value = "0.450 0.320 0.156 0.232 0.338 "
$thunderbird.Range("F51").Value = "0.769 0.843 0.913 0.922 0.977"

$thunderbird.Range("C52").Value = "0.679 0.605 0.282 0.434 0.460"
$thunderbird.Range("D52").Value = "0.477 0.487 0.303 0.339 0.414 "
$thunderbird.Range("E52").Value = "0.594 0.452 0.165 0.281 0.300 "
$thunderbird.Range("F52").Value = "0.744 0.839 0.912 0.906 0.971"

$thunderbird.Range("C53").Value = "0.595 0.641 0.395 0.460 0.475"
$thunderbird.Range("D53").Value = "0.555 0.493 0.263 0.385 0.556"
$thunderbird.Range("E53").Value = "0.448 0.496 0.250 0.303 0.313"
$thunderbird.Range("F53").Value = "0.779 0.841 0.895 0.913 0.976"

$thunderbird.Range("C54").Value = "0.402 0.245 0.044 0.060 0.400"
$thunderbird.Range("D54").Value = "0.405 0.456 0.125 0.280 0.556 "
$thunderbird.Range("E54").Value = "0.261 0.140 0.022 0.031 0.250"
$thunderbird.Range("F54").Value = "0.728 0.839 0.918 0.923 0.975"

$thunderbird.Range("C55").Value = "0.470 0.271 0.201 0.279 0.431"
$thunderbird.Range("D55").Value = "0.580 0.905 0.806 0.902 1.000"
$thunderbird.Range("E55").Value = "0.315 0.157 0.112 0.162 0.275"
$thunderbird.Range("F55").Value = "0.779 0.865 0.934 0.937 0.981"

$thunderbird.Range("C57").Value = "0.598 0.462 0.256 0.322 0.490 "
$thunderbird.Range("D57").Value = "0.524 0.517 0.277 0.324 0.722 "
$thunderbird.Range("E57").Value = "0.456 0.306 0.147 0.193 0.325"
$thunderbird.Range("F57").Value = "0.768 0.846 0.910 0.911 0.979"

$thunderbird.Range("C58").Value = "0.676 0.530 0.200 0.322 0.367 "
$thunderbird.Range("D58").Value = "0.490 0.516 0.316 0.423 0.643 "
$thunderbird.Range("E58").Value = "0.581 0.370 0.112 0.193 0.225 "
$thunderbird.Range("F58").Value = "0.752 0.847 0.918 0.921 0.977"

$thunderbird.Range("C59").Value = "0.609 0.627 0.373 0.481 0.490"
$thunderbird.Range("D59").Value = "0.566 0.490 0.268 0.424 0.703"
$thunderbird.Range("E59").Value = "0.463 0.479 0.232 0.320 0.325"
$thunderbird.Range("F59").Value = "0.784 0.840 0.898 0.918 0.979"

$thunderbird.Range("C60").Value = "0.375 0.223 0.035 0.051 0.367"
$thunderbird.Range("D60").Value = "0.396 0.445 0.114 0.162 0.621"
$thunderbird.Range("E60").Value = "0.238 0.126 0.018 0.026 0.225"
$thunderbird.Range("F60").Value = "0.727 0.838 0.919 0.918 0.976"

$thunderbird.Range("C61").Value = "0.448 0.281 0.215 0.279 0.431"
$thunderbird.Range("D61").Value = "0.642 0.940 0.871 0.925 1.000"
$thunderbird.Range("E61").Value = "0.293 0.163 0.121 0.162 0.275 "
$thunderbird.Range("F61").Value = "0.789 0.867 0.935 0.937 0.981"

# Commit also left the workbook with the Thunderbird tab selected/active
# (was Lucene before), scrolled/selected to the newly-added data.
$thunderbird.Activate()
$thunderbird.Range("D61").Select()
